$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous contents entirely (old range was A2:B11)
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "data length"
$ws.Range("B1").Value = "data list"

# Data row 2: A2 numeric, B2 string "A"
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "A"

# Column B list values 1..8 in rows 3..10
$values = 1,2,3,4,5,6,7,8
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
